# Updated cryptos list on Sat Oct 19 02:52:44 UTC 2024 with GitHub Actions
# Refreshes coin Price (D) and Volume(1h) (E) columns; ImmutableX/FirstDigitalUSD rows (36/37) swapped order.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.424.25"
$ws.Range("E2").Value = "  +0.60%  "

$ws.Range("D3").Value = "2.643.63"
$ws.Range("E3").Value = "  +0.42%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").Value = "'600.34"
$ws.Range("E5").Value = "  +0.65%  "

$ws.Range("D6").Value = "'154.71"
$ws.Range("E6").Value = "  +1.12%  "

$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("D8").Value = "'0.548"
$ws.Range("E8").Value = "  +0.78%  "

$ws.Range("D9").Value = "2.642.18"
$ws.Range("E9").Value = "  +0.44%  "

$ws.Range("E10").Value = "  +7.15%  "

$ws.Range("D11").Value = "'0.159"
$ws.Range("E11").Value = "  -0.40%  "

$ws.Range("D12").Value = "'5.26"
$ws.Range("E12").Value = "  +1.15%  "

$ws.Range("D13").Value = "'0.355"
$ws.Range("E13").Value = "  +1.49%  "

$ws.Range("D14").Value = "'28.27"
$ws.Range("E14").Value = "  +2.09%  "

$ws.Range("D15").Value = "'0.0000193"
$ws.Range("E15").Value = "  +0.58%  "

$ws.Range("D16").Value = "3.123.51"
$ws.Range("E16").Value = "  +0.31%  "

$ws.Range("D17").Value = "68.294.62"
$ws.Range("E17").Value = "  +0.81%  "

$ws.Range("D18").Value = "2.635.08"
$ws.Range("E18").Value = "  +0.35%  "

$ws.Range("D19").Value = "'11.51"
$ws.Range("E19").Value = "  +2.27%  "

$ws.Range("D20").Value = "'366.54"
$ws.Range("E20").Value = "  -1.71%  "

$ws.Range("D21").Value = "'7.54"
$ws.Range("E21").Value = "  +1.29%  "

$ws.Range("D22").Value = "'4.44"
$ws.Range("E22").Value = "  +4.49%  "

$ws.Range("D23").Value = "'4.89"
$ws.Range("E23").Value = "  +1.04%  "

$ws.Range("D24").Value = "'2.09"
$ws.Range("E24").Value = "  +1.44%  "

$ws.Range("D25").Value = "'73.96"
$ws.Range("E25").Value = "  +0.91%  "

$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  +0.04%  "

$ws.Range("D27").Value = "'9.92"
$ws.Range("E27").Value = "  +0.47%  "

$ws.Range("D28").Value = "'0.0000107"
$ws.Range("E28").Value = "  +1.77%  "

$ws.Range("D29").Value = "2.767.74"
$ws.Range("E29").Value = "  +0.15%  "

$ws.Range("D30").Value = "'1.00"
$ws.Range("E30").Value = "  +0.60%  "

$ws.Range("D31").Value = "'578.91"
$ws.Range("E31").Value = "  -2.12%  "

$ws.Range("D32").Value = "'8.12"
$ws.Range("E32").Value = "  +3.93%  "

$ws.Range("D33").Value = "'1.44"
$ws.Range("E33").Value = "  +4.11%  "

$ws.Range("D34").Value = "'1.88"
$ws.Range("E34").Value = "  +1.46%  "

$ws.Range("E35").Value = "  +3.92%  "

$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "'1.62"
$ws.Range("E36").Value = "  +5.81%  "

$ws.Range("B37").Value = "FirstDigitalUSD"
$ws.Range("C37").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D37").Value = "'0.999"
$ws.Range("E37").Value = "  -0.02%  "

$ws.Range("D38").Value = "'160.69"
$ws.Range("E38").Value = "  +1.22%  "

$ws.Range("D39").Value = "'19.51"
$ws.Range("E39").Value = "  +1.56%  "

$ws.Range("D40").Value = "'1.90"
$ws.Range("E40").Value = "  -0.20%  "

$ws.Range("D41").Value = "'0.373"
$ws.Range("E41").Value = "  +1.12%  "

$ws.Range("D42").Value = "'5.47"
$ws.Range("E42").Value = "  +3.20%  "

$ws.Range("D43").Value = "'2.69"
$ws.Range("E43").Value = "  +0.02%  "

$ws.Range("D44").Value = "'17.73"
$ws.Range("E44").Value = "  +3.42%  "

$ws.Range("D45").Value = "0.0₆0328"
$ws.Range("E45").Value = "  +10.32%  "

$ws.Range("E46").Value = "  -0.05%  "

$ws.Range("D47").Value = "'40.46"
$ws.Range("E47").Value = "  +0.12%  "

$ws.Range("D48").Value = "'157.74"
$ws.Range("E48").Value = "  +0.60%  "

$ws.Range("D49").Value = "'3.80"
$ws.Range("E49").Value = "  +2.73%  "

$ws.Range("D50").Value = "'1.72"
$ws.Range("E50").Value = "  +1.42%  "

$ws.Range("D51").Value = "'22.06"
$ws.Range("E51").Value = "  +3.32%  "
